# MLK_PMT_10104_-_V-004 : expand the parts table from 5 parts (rows 8-12)
# to 19 parts (rows 8-26), normalise text casing/units, and fix the
# equipment description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Equipment description fix -------------------------------------------
$ws.Range("D8").Value = "Air Receiver"

# --- Full target data for rows 8-26 (PARTS / SPEC / GRADE) ---------------
# columns: E=PARTS, I=SPEC., J=GRADE -- the rest (G,H,K,L,M,N,O) share the
# same value on every row of the table.
$parts = @(
    @{ E = "Shell Plate";       I = "ASTM A240"; J = "304L"   },
    @{ E = "Ellipsoidal Head";  I = "ASTM A240"; J = "304L"   },
    @{ E = "Flange (J)";        I = "ASTM A182"; J = "F304L"  },
    @{ E = "Flange (J2)";       I = "ASTM A182"; J = "F304L"  },
    @{ E = "Neck (K)";          I = "ASTM A312"; J = "TP304L" },
    @{ E = "Neck (L)";          I = "ASTM A312"; J = "TP304L" },
    @{ E = "Flange (M)";        I = "ASTM A182"; J = "F304L"  },
    @{ E = "Flange (M2)";       I = "ASTM A182"; J = "F304L"  },
    @{ E = "Neck (N)";          I = "ASTM A312"; J = "TP304L" },
    @{ E = "Flange (P)";        I = "ASTM A182"; J = "F304L"  },
    @{ E = "Neck (Q)";          I = "ASTM A312"; J = "TP304L" },
    @{ E = "Neck (R)";          I = "ASTM A312"; J = "TP304L" },
    @{ E = "Flange (S)";        I = "ASTM A182"; J = "F304L"  },
    @{ E = "Flange (S2)";       I = "ASTM A182"; J = "F304L"  },
    @{ E = "Neck (T)";          I = "ASTM A312"; J = "TP304L" },
    @{ E = "Neck (U)";          I = "ASTM A312"; J = "TP304L" },
    @{ E = "Blind Flange (V)";  I = "ASTM A182"; J = "F304L"  },
    @{ E = "Flange (X)";        I = "ASTM A182"; J = "F304L"  },
    @{ E = "Neck (Y)";          I = "ASTM A240"; J = "304L"   }
)

$firstRow = 8
for ($idx = 0; $idx -lt $parts.Count; $idx++) {
    $row = $firstRow + $idx
    $p = $parts[$idx]

    $ws.Cells.Item($row, 5).Value = $p.E                     # E - PARTS
    $ws.Cells.Item($row, 7).Value = "HOT WATER"               # G - FLUID
    $ws.Cells.Item($row, 8).Value = "Stainless Steel"         # H - MATERIAL
    $ws.Cells.Item($row, 9).Value = $p.I                      # I - SPEC.
    $ws.Cells.Item($row, 10).Value = $p.J                     # J - GRADE
    $ws.Cells.Item($row, 11).Value = "HOT INSULATION 100mm"   # K - INSULATION
    $ws.Cells.Item($row, 12).Value = "120 °C"                 # L - DESIGN TEMP
    $ws.Cells.Item($row, 13).Value = "4 Bar G"                # M - DESIGN PRESSURE
    $ws.Cells.Item($row, 14).Value = "90 °C"                  # N - OPERATING TEMP
    $ws.Cells.Item($row, 15).Value = "1 Bar G"                # O - OPERATING PRESSURE
}

# --- Re-merge the now-19-row-tall identity columns -------------------------
$ws.Range("A8:A26").Merge() | Out-Null
$ws.Range("B8:B26").Merge() | Out-Null
$ws.Range("C8:C26").Merge() | Out-Null
$ws.Range("D8:D26").Merge() | Out-Null
